# Update the player's match log so that it contains only a single match
# row (the "Oct 10 2020 vs Chennai Super Kings" entry), replacing the
# previous first data row, and remove the other data rows (rows 3-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every updated cell to be stored as text (matches the workbook's
# original numberStoredAsText convention) so numeric-looking values such
# as "22" or "157.14" are not reinterpreted/rounded as real numbers.
$row2 = $ws.Range("A2:K2")
$row2.NumberFormat = "@"

# Overwrite row 2 with the data that used to live in row 8, but keep the
# strike-rate value that was already present on row 2 (157.14).
$ws.Range("A2").Value = " Oct 10 2020"
$ws.Range("B2").Value = " Dubai (DSC)"
$ws.Range("C2").Value = "RCB won by 37 runs"
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Shivam Dube "
$ws.Range("G2").Value = "22"
$ws.Range("H2").Value = "14"
$ws.Range("I2").Value = "2"
$ws.Range("J2").Value = "1"
$ws.Range("K2").Value = "157.14"

# Remove the now-duplicated/obsolete rows 3 through 8.
$ws.Range("A3:K8").EntireRow.Delete()
